$wb = $excel.ActiveWorkbook
$wsCases = $wb.Worksheets.Item("Test_Cases")
$wsData  = $wb.Worksheets.Item("Test_Data")

# --- Test_Cases sheet: add "Functionality" column (E) values for rows 2-11 ---
$functionality = @(
    "Functionality_1",
    "Functionality_2",
    "Functionality_3",
    "Functionality_4",
    "Functionality_5",
    "Functionality_6",
    "Functionality_7",
    "Functionality_8",
    "Functionality_9",
    "Functionality_10"
)
for ($i = 0; $i -lt $functionality.Length; $i++) {
    $row = 2 + $i
    $wsCases.Cells.Item($row, 5).Value = $functionality[$i]
}

# Match the formatting of column D (border + centered) on the new E column cells
$wsCases.Range("D2:D11").Copy()
$wsCases.Range("E2:E11").PasteSpecial(-4122)

# Make the E1 header cell match the other header cells' formatting (A1:D1)
$wsCases.Range("A1").Copy()
$wsCases.Range("E1").PasteSpecial(-4122)

# Widen column E to match column C
$wsCases.Columns.Item(5).ColumnWidth = $wsCases.Columns.Item(3).ColumnWidth

# Highlight RegisterUser5-9 rows (B7:B11) with a light-orange fill
# (equivalent to Theme Color "Orange, Accent 2, Lighter 80%" -> RGB FBE5D6)
$wsCases.Range("B7:B11").Interior.Color = 14083579

# --- Restore / update selections on both sheets ---
[void]$wsData.Range("A17").Select()
[void]$wsCases.Range("C16").Select()
